$d = $word.ActiveDocument

$replacements = @(
    @("2023-08-17 Thursday", "2023-08-18 Friday"),
    @("95×30=2850", "95×81=7695"),
    @("87×81=7047", "37×78=2886"),
    @("30×29=870", "65×66=4290"),
    @("89×29=2581", "46×40=1840"),
    @("86×21=1806", "15×98=1470"),
    @("76×55=4180", "19×51=969"),
    @("72×78=5616", "58×53=3074"),
    @("46×21=966", "64×59=3776"),
    @("15×69=1035", "74×71=5254"),
    @("76×67=5092", "82×45=3690"),
    @("60×52=3120", "33×66=2178"),
    @("13×98=1274", "61×86=5246"),
    @("77×41=3157", "25×13=325"),
    @("90×50=4500", "80×59=4720"),
    @("14×17=238", "79×44=3476"),
    @("93×91=8463", "80×86=6880"),
    @("15×82=1230", "56×52=2912"),
    @("17×15=255", "77×37=2849"),
    @("34×22=748", "14×92=1288"),
    @("32×37=1184", "87×16=1392"),
    @("26×74=1924", "28×20=560"),
    @("23×63=1449", "25×55=1375"),
    @("90×97=8730", "19×30=570"),
    @("71×44=3124", "95×12=1140"),
    @("49×76=3724", "72×85=6120")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

$d.Save()
